# Update the "想去人数" (want-to-go count) values for a handful of events
# on the 展览 and 全部类型 sheets, matching the regenerated site data.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F7").Value  = 2508
$wsExpo.Range("F11").Value = 1510
$wsExpo.Range("F22").Value = 155
$wsExpo.Range("F23").Value = 48
$wsExpo.Range("F24").Value = 1591
$wsExpo.Range("F28").Value = 203

# Sheet "全部类型" (all types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F8").Value  = 2508
$wsAll.Range("F12").Value = 1510
$wsAll.Range("F23").Value = 155
$wsAll.Range("F24").Value = 48
$wsAll.Range("F25").Value = 1591
$wsAll.Range("F29").Value = 203

$wb.Save()
